{"js": "// The M2Doc test fixture embeds a Java stack trace (as literal text,\n// including \"\\t\" and \"\\n\" characters inside a single run) produced by a\n// test run. Re-running the test suite against a newer M2Doc build shifted\n// several line numbers in the M2Doc source files, changed a couple of\n// surefire/test-harness line numbers, and replaced the JDT/Eclipse test\n// runner tail of the trace with the Maven Surefire / Tycho / Equinox\n// launch chain. Apply each of those textual substitutions in place.\n\n// 1) Simple, unambiguous single-token replacements (line numbers / method\n//    accessor index) - each occurs exactly once, except the \":1096)\" one\n//    which occurs 3 times with the exact same old/new value everywhere.\nconst simpleReplacements = [\n  [\"M2DocEvaluator.java:1556)\", \"M2DocEvaluator.java:1703)\"],\n  [\"M2DocEvaluator.java:1096)\", \"M2DocEvaluator.java:1216)\"],\n  [\"M2DocEvaluator.java:1305)\", \"M2DocEvaluator.java:1425)\"],\n  [\"M2DocEvaluator.java:283)\", \"M2DocEvaluator.java:287)\"],\n  [\"M2DocEvaluator.java:272)\", \"M2DocEvaluator.java:276)\"],\n  [\"AbstractTemplatesTestSuite.java:479)\", \"AbstractTemplatesTestSuite.java:480)\"],\n  [\"AbstractTemplatesTestSuite.java:388)\", \"AbstractTemplatesTestSuite.java:389)\"],\n  [\"GeneratedMethodAccessor75.invoke\", \"GeneratedMethodAccessor74.invoke\"],\n];\n\nfor (const [oldText, newText] of simpleReplacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 2) The six-line JDT/Eclipse test-runner tail of the stack trace is\n//    replaced by a much longer Maven Surefire / Tycho / Equinox launch\n//    chain (the test now runs under a different harness).\nconst oldTail =\n  \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\\n\" +\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\\n\";\n\nconst newTail =\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)\\n\" +\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)\\n\" +\n  \"\\tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)\\n\" +\n  \"\\tat org.apache.maven.surefire.booter.ProviderFactory$ProviderProxy.invoke(ProviderFactory.java:156)\\n\" +\n  \"\\tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)\\n\" +\n  \"\\tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)\\n\" +\n  \"\\tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)\\n\" +\n  \"\\tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)\\n\" +\n  \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)\\n\" +\n  \"\\tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)\\n\" +\n  \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)\\n\" +\n  \"\\tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\\n\" +\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\\n\" +\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\\n\" +\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)\\n\" +\n  \"\\tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)\\n\";\n\nconst tailResults = context.document.body.search(oldTail, { matchCase: true });\ntailResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < tailResults.items.length; i++) {\n  tailResults.items[i].insertText(newTail, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The M2Doc test fixture embeds a Java stack trace (as literal text,\n# including tab/newline characters inside a single run) produced by a\n# test run. Re-running the test suite against a newer M2Doc build shifted\n# several line numbers in the M2Doc source files, changed a couple of\n# surefire/test-harness line numbers, and replaced the JDT/Eclipse test\n# runner tail of the trace with the Maven Surefire / Tycho / Equinox\n# launch chain. Apply each of those textual substitutions in place.\n\n$d = $word.ActiveDocument\n\n# 1) Simple, unambiguous single-token replacements (line numbers / method\n#    accessor index) - each occurs exactly once, except the \":1096)\" one\n#    which occurs 3 times with the exact same old/new value everywhere.\n$simpleReplacements = @(\n    @(\"M2DocEvaluator.java:1556)\", \"M2DocEvaluator.java:1703)\"),\n    @(\"M2DocEvaluator.java:1096)\", \"M2DocEvaluator.java:1216)\"),\n    @(\"M2DocEvaluator.java:1305)\", \"M2DocEvaluator.java:1425)\"),\n    @(\"M2DocEvaluator.java:283)\", \"M2DocEvaluator.java:287)\"),\n    @(\"M2DocEvaluator.java:272)\", \"M2DocEvaluator.java:276)\"),\n    @(\"AbstractTemplatesTestSuite.java:479)\", \"AbstractTemplatesTestSuite.java:480)\"),\n    @(\"AbstractTemplatesTestSuite.java:388)\", \"AbstractTemplatesTestSuite.java:389)\"),\n    @(\"GeneratedMethodAccessor75.invoke\", \"GeneratedMethodAccessor74.invoke\")\n)\n\nforeach ($pair in $simpleReplacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n\n# 2) The six-line JDT/Eclipse test-runner tail of the stack trace is\n#    replaced by a much longer Maven Surefire / Tycho / Equinox launch\n#    chain (the test now runs under a different harness).\n$oldTail = \"`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)`n\" + `\n    \"`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)`n\"\n\n$newTail = \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n\" + `\n    \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n\" + `\n    \"`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n    \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n    \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n    \"`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n\" + `\n    \"`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n\" + `\n    \"`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n\" + `\n    \"`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n\" + `\n    \"`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n    \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n    \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n    \"`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n\" + `\n    \"`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n\" + `\n    \"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n\" + `\n    \"`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n\" + `\n    \"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n\" + `\n    \"`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n\" + `\n    \"`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n\" + `\n    \"`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n\" + `\n    \"`tat java.lang.reflect.Method.invoke(Method.java:498)`n\" + `\n    \"`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n\" + `\n    \"`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n\" + `\n    \"`tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n\" + `\n    \"`tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)`n\"\n\n$find = $d.Content.Find\n$find.Text = $oldTail\n$find.Replacement.Text = $newTail\n$find.Execute([ref]$oldTail, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newTail, 2) | Out-Null\n"}
